$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values to match refreshed cryptocurrency data.
# NumberFormat "@" (Text) is applied before each write and the style is
# reset to "Normal" right after, so numeric-looking strings (prices like
# "1.001" or "238.08") are stored as literal text -- matching the source
# feed -- instead of being auto-coerced into Excel numbers, while leaving
# no residual formatting/style on the cell.

# Row 2
$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '29.130.90'
$cell.Style = "Normal"
$cell = $ws.Range('E2')
$cell.NumberFormat = "@"
$cell.Value = '  -3.28%  '
$cell.Style = "Normal"
# Row 3
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '1.849.59'
$cell.Style = "Normal"
$cell = $ws.Range('E3')
$cell.NumberFormat = "@"
$cell.Value = '  -2.17%  '
$cell.Style = "Normal"
# Row 4
$cell = $ws.Range('D4')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$cell = $ws.Range('E4')
$cell.NumberFormat = "@"
$cell.Value = '  -0.05%  '
$cell.Style = "Normal"
# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '0.7032'
$cell.Style = "Normal"
$cell = $ws.Range('E5')
$cell.NumberFormat = "@"
$cell.Value = '  -5.00%  '
$cell.Style = "Normal"
# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '238.08'
$cell.Style = "Normal"
$cell = $ws.Range('E6')
$cell.NumberFormat = "@"
$cell.Value = '  -1.82%  '
$cell.Style = "Normal"
# Row 7
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$cell = $ws.Range('E7')
$cell.NumberFormat = "@"
$cell.Value = '  -0.11%  '
$cell.Style = "Normal"
# Row 8
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.3036'
$cell.Style = "Normal"
$cell = $ws.Range('E8')
$cell.NumberFormat = "@"
$cell.Value = '  -4.09%  '
$cell.Style = "Normal"
# Row 9
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.07474'
$cell.Style = "Normal"
$cell = $ws.Range('E9')
$cell.NumberFormat = "@"
$cell.Value = '  +3.74%  '
$cell.Style = "Normal"
# Row 10
$cell = $ws.Range('E10')
$cell.NumberFormat = "@"
$cell.Value = '  -6.44%  '
$cell.Style = "Normal"
# Row 11
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.08122'
$cell.Style = "Normal"
$cell = $ws.Range('E11')
$cell.NumberFormat = "@"
$cell.Value = '  -2.73%  '
$cell.Style = "Normal"
# Row 12
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '0.7246'
$cell.Style = "Normal"
$cell = $ws.Range('E12')
$cell.NumberFormat = "@"
$cell.Value = '  -4.67%  '
$cell.Style = "Normal"
# Row 13
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '1.829.84'
$cell.Style = "Normal"
$cell = $ws.Range('E13')
$cell.NumberFormat = "@"
$cell.Value = '  -5.04%  '
$cell.Style = "Normal"
# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '5.210'
$cell.Style = "Normal"
$cell = $ws.Range('E14')
$cell.NumberFormat = "@"
$cell.Value = '  -4.12%  '
$cell.Style = "Normal"
# Row 15
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '88.91'
$cell.Style = "Normal"
$cell = $ws.Range('E15')
$cell.NumberFormat = "@"
$cell.Value = '  -4.19%  '
$cell.Style = "Normal"
# Row 16
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '29.046.35'
$cell.Style = "Normal"
# Row 17
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '5.769'
$cell.Style = "Normal"
$cell = $ws.Range('E17')
$cell.NumberFormat = "@"
$cell.Value = '  -6.72%  '
$cell.Style = "Normal"
# Row 18
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '238.36'
$cell.Style = "Normal"
$cell = $ws.Range('E18')
$cell.NumberFormat = "@"
$cell.Value = '  -4.68%  '
$cell.Style = "Normal"
# Row 19
$cell = $ws.Range('B19')
$cell.NumberFormat = "@"
$cell.Value = 'ShibaInu'
$cell.Style = "Normal"
$cell = $ws.Range('C19')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell.Style = "Normal"
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '0.000007661'
$cell.Style = "Normal"
$cell = $ws.Range('E19')
$cell.NumberFormat = "@"
$cell.Value = '  -2.56%  '
$cell.Style = "Normal"
# Row 20
$cell = $ws.Range('B20')
$cell.NumberFormat = "@"
$cell.Value = 'Avalanche'
$cell.Style = "Normal"
$cell = $ws.Range('C20')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell.Style = "Normal"
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '13.05'
$cell.Style = "Normal"
$cell = $ws.Range('E20')
$cell.NumberFormat = "@"
$cell.Value = '  -4.36%  '
$cell.Style = "Normal"
# Row 21
$cell = $ws.Range('E21')
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell.Style = "Normal"
# Row 22
$cell = $ws.Range('B22')
$cell.NumberFormat = "@"
$cell.Value = 'BinanceUSD'
$cell.Style = "Normal"
$cell = $ws.Range('C22')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell.Style = "Normal"
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$cell = $ws.Range('E22')
$cell.NumberFormat = "@"
$cell.Value = '  -0.02%  '
$cell.Style = "Normal"
# Row 23
$cell = $ws.Range('B23')
$cell.NumberFormat = "@"
$cell.Value = 'WrappedliquidstakedEther2.0'
$cell.Style = "Normal"
$cell = $ws.Range('C23')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$cell.Style = "Normal"
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '2.076.95'
$cell.Style = "Normal"
$cell = $ws.Range('E23')
$cell.NumberFormat = "@"
$cell.Value = '  -5.79%  '
$cell.Style = "Normal"
# Row 24
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '7.539'
$cell.Style = "Normal"
$cell = $ws.Range('E24')
$cell.NumberFormat = "@"
$cell.Value = '  -5.74%  '
$cell.Style = "Normal"
# Row 25
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '161.81'
$cell.Style = "Normal"
$cell = $ws.Range('E25')
$cell.NumberFormat = "@"
$cell.Value = '  -1.39%  '
$cell.Style = "Normal"
# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '8.977'
$cell.Style = "Normal"
$cell = $ws.Range('E26')
$cell.NumberFormat = "@"
$cell.Value = '  -3.41%  '
$cell.Style = "Normal"
# Row 27
$cell = $ws.Range('E27')
$cell.NumberFormat = "@"
$cell.Value = '  -7.60%  '
$cell.Style = "Normal"
# Row 28
$cell = $ws.Range('E28')
$cell.NumberFormat = "@"
$cell.Value = '  -4.00%  '
$cell.Style = "Normal"
# Row 29
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '1.938'
$cell.Style = "Normal"
$cell = $ws.Range('E29')
$cell.NumberFormat = "@"
$cell.Value = '  -6.04%  '
$cell.Style = "Normal"
# Row 30
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '1.385'
$cell.Style = "Normal"
# Row 31
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '4.543'
$cell.Style = "Normal"
$cell = $ws.Range('E31')
$cell.NumberFormat = "@"
$cell.Value = '  -1.20%  '
$cell.Style = "Normal"
# Row 32
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '1.494'
$cell.Style = "Normal"
$cell = $ws.Range('E32')
$cell.NumberFormat = "@"
$cell.Value = '  -2.65%  '
$cell.Style = "Normal"
# Row 33
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '3.987'
$cell.Style = "Normal"
$cell = $ws.Range('E33')
$cell.NumberFormat = "@"
$cell.Value = '  -5.45%  '
$cell.Style = "Normal"
# Row 34
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '0.05148'
$cell.Style = "Normal"
$cell = $ws.Range('E34')
$cell.NumberFormat = "@"
$cell.Value = '  -4.57%  '
$cell.Style = "Normal"
# Row 35
$cell = $ws.Range('E35')
$cell.NumberFormat = "@"
$cell.Value = '  -5.15%  '
$cell.Style = "Normal"
# Row 36
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '1.037'
$cell.Style = "Normal"
$cell = $ws.Range('E36')
$cell.NumberFormat = "@"
$cell.Value = '  +3.53%  '
$cell.Style = "Normal"
# Row 37
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '0.7002'
$cell.Style = "Normal"
$cell = $ws.Range('E37')
$cell.NumberFormat = "@"
$cell.Value = '  -8.64%  '
$cell.Style = "Normal"
# Row 38
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '2.642'
$cell.Style = "Normal"
$cell = $ws.Range('E38')
$cell.NumberFormat = "@"
$cell.Value = '  -3.00%  '
$cell.Style = "Normal"
# Row 39
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.01868'
$cell.Style = "Normal"
$cell = $ws.Range('E39')
$cell.NumberFormat = "@"
$cell.Value = '  -5.17%  '
$cell.Style = "Normal"
# Row 40
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '2.679'
$cell.Style = "Normal"
$cell = $ws.Range('E40')
$cell.NumberFormat = "@"
$cell.Value = '  -3.11%  '
$cell.Style = "Normal"
# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.9557'
$cell.Style = "Normal"
$cell = $ws.Range('E41')
$cell.NumberFormat = "@"
$cell.Value = '  +9.90%  '
$cell.Style = "Normal"
# Row 42
$cell = $ws.Range('B42')
$cell.NumberFormat = "@"
$cell.Value = 'Maker'
$cell.Style = "Normal"
$cell = $ws.Range('C42')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell.Style = "Normal"
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '1.085.37'
$cell.Style = "Normal"
$cell = $ws.Range('E42')
$cell.NumberFormat = "@"
$cell.Value = '  -1.64%  '
$cell.Style = "Normal"
# Row 43
$cell = $ws.Range('B43')
$cell.NumberFormat = "@"
$cell.Value = 'FraxShare'
$cell.Style = "Normal"
$cell = $ws.Range('C43')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell.Style = "Normal"
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '5.994'
$cell.Style = "Normal"
$cell = $ws.Range('E43')
$cell.NumberFormat = "@"
$cell.Value = '  -1.27%  '
$cell.Style = "Normal"
# Row 44
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '0.4289'
$cell.Style = "Normal"
$cell = $ws.Range('E44')
$cell.NumberFormat = "@"
$cell.Value = '  -5.75%  '
$cell.Style = "Normal"
# Row 45
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '69.81'
$cell.Style = "Normal"
$cell = $ws.Range('E45')
$cell.NumberFormat = "@"
$cell.Value = '  -3.82%  '
$cell.Style = "Normal"
# Row 46
$cell = $ws.Range('E46')
$cell.NumberFormat = "@"
$cell.Value = '  -0.14%  '
$cell.Style = "Normal"
# Row 47
$cell = $ws.Range('E47')
$cell.NumberFormat = "@"
$cell.Value = '  -2.08%  '
$cell.Style = "Normal"
# Row 48
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '1.742'
$cell.Style = "Normal"
$cell = $ws.Range('E48')
$cell.NumberFormat = "@"
$cell.Value = '  -6.62%  '
$cell.Style = "Normal"
# Row 49
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '1.980.49'
$cell.Style = "Normal"
$cell = $ws.Range('E49')
$cell.NumberFormat = "@"
$cell.Value = '  -4.11%  '
$cell.Style = "Normal"
# Row 50
$cell = $ws.Range('B50')
$cell.NumberFormat = "@"
$cell.Value = 'Aptos'
$cell.Style = "Normal"
$cell = $ws.Range('C50')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell.Style = "Normal"
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '7.040'
$cell.Style = "Normal"
$cell = $ws.Range('E50')
$cell.NumberFormat = "@"
$cell.Value = '  -7.40%  '
$cell.Style = "Normal"
# Row 51
$cell = $ws.Range('B51')
$cell.NumberFormat = "@"
$cell.Value = 'EnergySwap'
$cell.Style = "Normal"
$cell = $ws.Range('C51')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell.Style = "Normal"
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '9.142'
$cell.Style = "Normal"
$cell = $ws.Range('E51')
$cell.NumberFormat = "@"
$cell.Value = '  -4.97%  '
$cell.Style = "Normal"
